$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "2" element rows (and the one-off "charging_station1"
# row) that are no longer produced by the chart-creating function. Delete
# from the bottom up so earlier row numbers stay valid while we work.
$rowsToDelete = @(21, 17, 15, 14, 13, 11, 9, 5)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).EntireRow.Delete()
}
